$d = $word.ActiveDocument

# Update the date in the title paragraph
$d.Content.Find.Execute("2025-08-21 Thursday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-08-22 Friday", 2) | Out-Null

# Update each table cell by explicit row/column address to avoid ambiguous text collisions
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "54+3="
$t.Cell(1, 2).Range.Text = "38-11="
$t.Cell(1, 3).Range.Text = "94-72="
$t.Cell(1, 4).Range.Text = "10+69="
$t.Cell(1, 5).Range.Text = "45-7="
$t.Cell(2, 1).Range.Text = "79+5="
$t.Cell(2, 2).Range.Text = "91-28="
$t.Cell(2, 3).Range.Text = "8+76="
$t.Cell(2, 4).Range.Text = "86-52="
$t.Cell(2, 5).Range.Text = "21+24="
$t.Cell(3, 1).Range.Text = "74-72="
$t.Cell(3, 2).Range.Text = "19+62="
$t.Cell(3, 3).Range.Text = "63+32="
$t.Cell(3, 4).Range.Text = "23+71="
$t.Cell(3, 5).Range.Text = "9+57="
$t.Cell(4, 1).Range.Text = "44+14="
$t.Cell(4, 2).Range.Text = "61+35="
$t.Cell(4, 3).Range.Text = "45+48="
$t.Cell(4, 4).Range.Text = "23+14="
$t.Cell(4, 5).Range.Text = "6+75="
$t.Cell(5, 1).Range.Text = "72-58="
$t.Cell(5, 2).Range.Text = "32-8="
$t.Cell(5, 3).Range.Text = "44+44="
$t.Cell(5, 4).Range.Text = "14+80="
$t.Cell(5, 5).Range.Text = "16+75="
$t.Cell(6, 1).Range.Text = "97-1="
$t.Cell(6, 2).Range.Text = "75-15="
$t.Cell(6, 3).Range.Text = "97-92="
$t.Cell(6, 4).Range.Text = "23-1="
$t.Cell(6, 5).Range.Text = "0+22="
$t.Cell(7, 1).Range.Text = "82-70="
$t.Cell(7, 2).Range.Text = "76-53="
$t.Cell(7, 3).Range.Text = "36+39="
$t.Cell(7, 4).Range.Text = "12+10="
$t.Cell(7, 5).Range.Text = "36+17="
$t.Cell(8, 1).Range.Text = "39+59="
$t.Cell(8, 2).Range.Text = "46+34="
$t.Cell(8, 3).Range.Text = "4+38="
$t.Cell(8, 4).Range.Text = "91-20="
$t.Cell(8, 5).Range.Text = "72+25="
$t.Cell(9, 1).Range.Text = "76-35="
$t.Cell(9, 2).Range.Text = "97-68="
$t.Cell(9, 3).Range.Text = "13+65="
$t.Cell(9, 4).Range.Text = "23+22="
$t.Cell(9, 5).Range.Text = "93-55="
$t.Cell(10, 1).Range.Text = "77+14="
$t.Cell(10, 2).Range.Text = "99-46="
$t.Cell(10, 3).Range.Text = "98-72="
$t.Cell(10, 4).Range.Text = "22-11="
$t.Cell(10, 5).Range.Text = "35-14="
$t.Cell(11, 1).Range.Text = "19-1="
$t.Cell(11, 2).Range.Text = "48-34="
$t.Cell(11, 3).Range.Text = "88-48="
$t.Cell(11, 4).Range.Text = "4+48="
$t.Cell(11, 5).Range.Text = "29-25="
$t.Cell(12, 1).Range.Text = "2+59="
$t.Cell(12, 2).Range.Text = "20+1="
$t.Cell(12, 3).Range.Text = "5-2="
$t.Cell(12, 4).Range.Text = "72-36="
$t.Cell(12, 5).Range.Text = "23-2="
$t.Cell(13, 1).Range.Text = "69-48="
$t.Cell(13, 2).Range.Text = "55-1="
$t.Cell(13, 3).Range.Text = "13+63="
$t.Cell(13, 4).Range.Text = "49-3="
$t.Cell(13, 5).Range.Text = "85-44="
$t.Cell(14, 1).Range.Text = "10+53="
$t.Cell(14, 2).Range.Text = "16+49="
$t.Cell(14, 3).Range.Text = "48-17="
$t.Cell(14, 4).Range.Text = "43+13="
$t.Cell(14, 5).Range.Text = "68-14="
$t.Cell(15, 1).Range.Text = "41+53="
$t.Cell(15, 2).Range.Text = "76+4="
$t.Cell(15, 3).Range.Text = "94+2="
$t.Cell(15, 4).Range.Text = "93-39="
$t.Cell(15, 5).Range.Text = "77+1="
$t.Cell(16, 1).Range.Text = "16+27="
$t.Cell(16, 2).Range.Text = "79+20="
$t.Cell(16, 3).Range.Text = "0+46="
$t.Cell(16, 4).Range.Text = "1+12="
$t.Cell(16, 5).Range.Text = "26-19="
$t.Cell(17, 1).Range.Text = "68-60="
$t.Cell(17, 2).Range.Text = "54-41="
$t.Cell(17, 3).Range.Text = "10+48="
$t.Cell(17, 4).Range.Text = "26+57="
$t.Cell(17, 5).Range.Text = "9+12="
$t.Cell(18, 1).Range.Text = "72+8="
$t.Cell(18, 2).Range.Text = "52-13="
$t.Cell(18, 3).Range.Text = "22+40="
$t.Cell(18, 4).Range.Text = "94-45="
$t.Cell(18, 5).Range.Text = "16+60="
$t.Cell(19, 1).Range.Text = "50-10="
$t.Cell(19, 2).Range.Text = "80-71="
$t.Cell(19, 3).Range.Text = "71-13="
$t.Cell(19, 4).Range.Text = "33-10="
$t.Cell(19, 5).Range.Text = "62-23="
$t.Cell(20, 1).Range.Text = "56-35="
$t.Cell(20, 2).Range.Text = "47-42="
$t.Cell(20, 3).Range.Text = "54-28="
$t.Cell(20, 4).Range.Text = "11+84="
$t.Cell(20, 5).Range.Text = "80-61="
